$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (A and B) ---
# COM ColumnWidth is expressed in "characters" (stored XML width = ColumnWidth + 5/MaxDigitWidth-ish offset),
# so back the stored/target width out by the standard 5/6 offset used by this engine.
$ws.Columns.Item(1).ColumnWidth = 39.42578125 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 38 - (5/6)

# --- New column S data (years 2022 row, with matching figures) ---
# Row 4 (header year row) - reuse formatting from R4
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

# Row 5 - reuse formatting from R5
$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 44

# Row 6 - reuse formatting from R6, then apply the "0.0" number format (new combined style)
$ws.Range("R6").Copy()
$ws.Range("S6").PasteSpecial(-4122)
$ws.Range("S6").NumberFormat = "0.0"
$ws.Range("S6").Value = 20.6

# Row 7 - reuse formatting from R7
$ws.Range("R7").Copy()
$ws.Range("S7").PasteSpecial(-4122)
$ws.Range("S7").Value = 7.9

# Row 8 - reuse formatting from R8
$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("S8").Value = 15.5

# --- Selection cursor moved by the editor ---
$ws.Range("Y14").Select()
